$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62: both Q and R "detect_structure"/"backup" flags reset to 0
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = 0

# Remaining rows: reset the "backup" (R) flag to 0 (previously 1 or 2)
$backupResetRows = @(69, 80, 91, 94, 99, 116, 128, 136, 147, 157, 160, 166, 174, 183, 197, 208, 218, 226, 241, 250, 253, 262, 266, 275, 289, 307, 318, 320, 329, 341, 344, 352, 364, 384, 410, 430, 439, 440, 453, 478, 480, 488, 499, 515, 521, 535, 538, 548, 552, 557, 559, 565, 573, 579, 591, 598, 610, 627, 634, 637, 643, 646, 653, 660, 677, 687, 695, 709, 714, 728, 730, 739, 744, 753, 761, 783, 785, 791, 802, 807, 815, 823, 830, 837, 843, 849, 863, 870, 878, 895, 907, 914, 918, 923, 926, 935, 936, 944, 955, 965, 976, 985, 994, 1001, 1007, 1010, 1018, 1038, 1040, 1051, 1057, 1072, 1086, 1098, 1103, 1113, 1119, 1125, 1142, 1146, 1152, 1156, 1169, 1171, 1178, 1185, 1195, 1196, 1214, 1221, 1227, 1234, 1242, 1245, 1260, 1272, 1293, 1298, 1316, 1331, 1336, 1342, 1354, 1360, 1377, 1389, 1401, 1403, 1415, 1419, 1428, 1434, 1437, 1459, 1469, 1475)
foreach ($r in $backupResetRows) {
    $ws.Cells.Item($r, 18).Value = 0
}

# Row 1485: isPivot flag set (reclassified pivot week)
$ws.Range("O1485").Value = 3

# Rows 1487-1488: previously-blank "backup" flags get explicit 0 now that they are processed
$ws.Range("R1487").Value = 0
$ws.Range("R1488").Value = 0

# New weekly rows appended: 1489-1499 (2024-07-01 .. 2024-09-09)
# Column A uses the same custom datetime number format as the rest of the date column
$ws.Range("A1489:A1499").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 1489
$ws.Range("A1489").Value = 45474
$ws.Range("B1489").Value = 699
$ws.Range("C1489").Value = 705.75
$ws.Range("D1489").Value = 683.4000244140625
$ws.Range("E1489").Value = 699.1500244140625
$ws.Range("F1489").Value = 695.1649780273438
$ws.Range("G1489").Value = 27327302
$ws.Range("H1489").Value = 2024
$ws.Range("I1489").Value = 7
$ws.Range("J1489").Value = 1
$ws.Range("K1489").Value = 0
$ws.Range("L1489").Value = 0
$ws.Range("M1489").Value = 0
$ws.Range("N1489").Value = 27
$ws.Range("O1489").Value = 0
$ws.Range("P1489").Value = 0
$ws.Range("Q1489").Value = 0

# Row 1490
$ws.Range("A1490").Value = 45481
$ws.Range("B1490").Value = 703.6500244140625
$ws.Range("C1490").Value = 712.0999755859375
$ws.Range("D1490").Value = 685.2999877929688
$ws.Range("E1490").Value = 692.0499877929688
$ws.Range("F1490").Value = 688.1054077148438
$ws.Range("G1490").Value = 24167004
$ws.Range("H1490").Value = 2024
$ws.Range("I1490").Value = 7
$ws.Range("J1490").Value = 8
$ws.Range("K1490").Value = 0
$ws.Range("L1490").Value = 0
$ws.Range("M1490").Value = 0
$ws.Range("N1490").Value = 28
$ws.Range("O1490").Value = 0
$ws.Range("P1490").Value = 0
$ws.Range("Q1490").Value = 0

# Row 1491
$ws.Range("A1491").Value = 45488
$ws.Range("B1491").Value = 692.0499877929688
$ws.Range("C1491").Value = 703.1500244140625
$ws.Range("D1491").Value = 658.6500244140625
$ws.Range("E1491").Value = 663
$ws.Range("F1491").Value = 659.2210083007812
$ws.Range("G1491").Value = 21467327
$ws.Range("H1491").Value = 2024
$ws.Range("I1491").Value = 7
$ws.Range("J1491").Value = 15
$ws.Range("K1491").Value = 0
$ws.Range("L1491").Value = 0
$ws.Range("M1491").Value = 0
$ws.Range("N1491").Value = 29
$ws.Range("O1491").Value = 0
$ws.Range("P1491").Value = 0
$ws.Range("Q1491").Value = 0

# Row 1492
$ws.Range("A1492").Value = 45495
$ws.Range("B1492").Value = 655
$ws.Range("C1492").Value = 675
$ws.Range("D1492").Value = 633.2999877929688
$ws.Range("E1492").Value = 667.5999755859375
$ws.Range("F1492").Value = 663.7947387695312
$ws.Range("G1492").Value = 26406929
$ws.Range("H1492").Value = 2024
$ws.Range("I1492").Value = 7
$ws.Range("J1492").Value = 22
$ws.Range("K1492").Value = 0
$ws.Range("L1492").Value = 0
$ws.Range("M1492").Value = 0
$ws.Range("N1492").Value = 30
$ws.Range("O1492").Value = 0
$ws.Range("P1492").Value = 0
$ws.Range("Q1492").Value = 0

# Row 1493
$ws.Range("A1493").Value = 45502
$ws.Range("B1493").Value = 676.75
$ws.Range("C1493").Value = 690.9000244140625
$ws.Range("D1493").Value = 646.5499877929688
$ws.Range("E1493").Value = 648.0499877929688
$ws.Range("F1493").Value = 644.356201171875
$ws.Range("G1493").Value = 36847249
$ws.Range("H1493").Value = 2024
$ws.Range("I1493").Value = 7
$ws.Range("J1493").Value = 29
$ws.Range("K1493").Value = 0
$ws.Range("L1493").Value = 0
$ws.Range("M1493").Value = 0
$ws.Range("N1493").Value = 31
$ws.Range("O1493").Value = 0
$ws.Range("P1493").Value = 0
$ws.Range("Q1493").Value = 0

# Row 1494
$ws.Range("A1494").Value = 45509
$ws.Range("B1494").Value = 628.0499877929688
$ws.Range("C1494").Value = 633.8499755859375
$ws.Range("D1494").Value = 608
$ws.Range("E1494").Value = 622.9000244140625
$ws.Range("F1494").Value = 619.3495483398438
$ws.Range("G1494").Value = 34380122
$ws.Range("H1494").Value = 2024
$ws.Range("I1494").Value = 8
$ws.Range("J1494").Value = 5
$ws.Range("K1494").Value = 0
$ws.Range("L1494").Value = 0
$ws.Range("M1494").Value = 0
$ws.Range("N1494").Value = 32
$ws.Range("O1494").Value = 2
$ws.Range("P1494").Value = 0
$ws.Range("Q1494").Value = 0

# Row 1495
$ws.Range("A1495").Value = 45516
$ws.Range("B1495").Value = 619
$ws.Range("C1495").Value = 642
$ws.Range("D1495").Value = 611.5999755859375
$ws.Range("E1495").Value = 634.1500244140625
$ws.Range("F1495").Value = 634.1500244140625
$ws.Range("G1495").Value = 25592841
$ws.Range("H1495").Value = 2024
$ws.Range("I1495").Value = 8
$ws.Range("J1495").Value = 12
$ws.Range("K1495").Value = 0
$ws.Range("L1495").Value = 0
$ws.Range("M1495").Value = 0
$ws.Range("N1495").Value = 33
$ws.Range("O1495").Value = 0
$ws.Range("P1495").Value = 0
$ws.Range("Q1495").Value = 0

# Row 1496
$ws.Range("A1496").Value = 45523
$ws.Range("B1496").Value = 636.5
$ws.Range("C1496").Value = 696.5999755859375
$ws.Range("D1496").Value = 636.5
$ws.Range("E1496").Value = 685.0999755859375
$ws.Range("F1496").Value = 685.0999755859375
$ws.Range("G1496").Value = 41584996
$ws.Range("H1496").Value = 2024
$ws.Range("I1496").Value = 8
$ws.Range("J1496").Value = 19
$ws.Range("K1496").Value = 0
$ws.Range("L1496").Value = 0
$ws.Range("M1496").Value = 0
$ws.Range("N1496").Value = 34
$ws.Range("O1496").Value = 0
$ws.Range("P1496").Value = 0
$ws.Range("Q1496").Value = 0

# Row 1497
$ws.Range("A1497").Value = 45530
$ws.Range("B1497").Value = 694.4000244140625
$ws.Range("C1497").Value = 712.9500122070312
$ws.Range("D1497").Value = 689.5
$ws.Range("E1497").Value = 701.3499755859375
$ws.Range("F1497").Value = 701.3499755859375
$ws.Range("G1497").Value = 34715868
$ws.Range("H1497").Value = 2024
$ws.Range("I1497").Value = 8
$ws.Range("J1497").Value = 26
$ws.Range("K1497").Value = 0
$ws.Range("L1497").Value = 0
$ws.Range("M1497").Value = 0
$ws.Range("N1497").Value = 35
$ws.Range("O1497").Value = 0
$ws.Range("P1497").Value = 0
$ws.Range("Q1497").Value = 0

# Row 1498
$ws.Range("A1498").Value = 45537
$ws.Range("B1498").Value = 706.7000122070312
$ws.Range("C1498").Value = 706.7000122070312
$ws.Range("D1498").Value = 661.25
$ws.Range("E1498").Value = 667.0999755859375
$ws.Range("F1498").Value = 667.0999755859375
$ws.Range("G1498").Value = 26521523
$ws.Range("H1498").Value = 2024
$ws.Range("I1498").Value = 9
$ws.Range("J1498").Value = 2
$ws.Range("K1498").Value = 0
$ws.Range("L1498").Value = 0
$ws.Range("M1498").Value = 0
$ws.Range("N1498").Value = 36
$ws.Range("O1498").Value = 0
$ws.Range("P1498").Value = 0
$ws.Range("Q1498").Value = 0

# Row 1499
$ws.Range("A1499").Value = 45544
$ws.Range("B1499").Value = 664
$ws.Range("C1499").Value = 686.4000244140625
$ws.Range("D1499").Value = 645.4000244140625
$ws.Range("E1499").Value = 673.25
$ws.Range("F1499").Value = 673.25
$ws.Range("G1499").Value = 26844184
$ws.Range("H1499").Value = 2024
$ws.Range("I1499").Value = 9
$ws.Range("J1499").Value = 9
$ws.Range("K1499").Value = 0
$ws.Range("L1499").Value = 0
$ws.Range("M1499").Value = 0
$ws.Range("N1499").Value = 37
$ws.Range("O1499").Value = 0
$ws.Range("P1499").Value = 0
$ws.Range("Q1499").Value = 0
